# edit.ps1 - weekly CompStat (33rd Precinct) refresh: new reporting week + updated crime counts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bulletin volume/number and the reporting week date range ---
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# --- Cells that flip between the "no data" text placeholder ("0" / "***.*") and a
#     real number need their number-format class changed too. Copy the value+format
#     from a stable donor cell that already carries the desired class, then (for
#     numeric targets) overwrite with the correct figure.
$ws.Range("G14").Copy($ws.Range("D14"))
$ws.Range("D14").Value = 1
$ws.Range("H14").Copy($ws.Range("E14"))
$ws.Range("E14").Value = -100
$ws.Range("G14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 2
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("N14").Copy($ws.Range("E15"))
$ws.Range("G14").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 2
$ws.Range("G14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("H14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("G14").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("H14").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -100
$ws.Range("G14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("N14").Copy($ws.Range("E27"))
$ws.Range("G14").Copy($ws.Range("F27"))
$ws.Range("F27").Value = 2
$ws.Range("G14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("C14").Copy($ws.Range("C30"))
$ws.Range("G14").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("H14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("G14").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("H14").Copy($ws.Range("H31"))
$ws.Range("H31").Value = -100

# --- Remaining cells: value changes only, number-format class unchanged ---
$ws.Range("J14").Value = 4
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 60
$ws.Range("L15").Value = 60
$ws.Range("M15").Value = -20
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 63
$ws.Range("J16").Value = 86
$ws.Range("K16").Value = -26.744186046511
$ws.Range("L16").Value = -3.076923076923
$ws.Range("M16").Value = -43.243243243243
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 27.777777777777
$ws.Range("I17").Value = 116
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = -3.333333333333
$ws.Range("L17").Value = 1.754385964912
$ws.Range("M17").Value = 41.463414634146
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = -13.725490196078
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -10.204081632653
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 30
$ws.Range("H19").Value = 11.111111111111
$ws.Range("I19").Value = 159
$ws.Range("J19").Value = 171
$ws.Range("K19").Value = -7.017543859649
$ws.Range("L19").Value = 7.432432432432
$ws.Range("M19").Value = 22.307692307692
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -57.142857142857
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 27.272727272727
$ws.Range("I20").Value = 64
$ws.Range("J20").Value = 51
$ws.Range("K20").Value = 25.490196078431
$ws.Range("L20").Value = -4.477611940298
$ws.Range("M20").Value = 100
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 20.833333333333
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = 22.784810126582
$ws.Range("I21").Value = 454
$ws.Range("J21").Value = 488
$ws.Range("K21").Value = -6.967213114754
$ws.Range("L21").Value = -2.783725910064
$ws.Range("M21").Value = 8.872901678657
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = -33.333333333333
$ws.Range("M22").Value = 0
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = -42.857142857142
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -52.941176470588
$ws.Range("F24").Value = 38
$ws.Range("H24").Value = -54.216867469879
$ws.Range("I24").Value = 360
$ws.Range("J24").Value = 436
$ws.Range("K24").Value = -17.431192660550
$ws.Range("L24").Value = -20.353982300885
$ws.Range("M24").Value = 44
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -83.333333333333
$ws.Range("I25").Value = 85
$ws.Range("J25").Value = 131
$ws.Range("K25").Value = -35.114503816793
$ws.Range("L25").Value = -25.438596491228
$ws.Range("C26").Value = 4
$ws.Range("E26").Value = -63.636363636363
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 176
$ws.Range("J26").Value = 162
$ws.Range("K26").Value = 8.641975308641
$ws.Range("L26").Value = -7.368421052631
$ws.Range("M26").Value = -22.807017543859
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 14.285714285714
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 29
$ws.Range("K28").Value = 61.111111111111
$ws.Range("L28").Value = 70.588235294117
$ws.Range("L29").Value = -85.714285714285
$ws.Range("L30").Value = -80
$ws.Range("J31").Value = 2
